# Auto-generated edit script applying the Kujata_Profits value update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1559
$ws.Range("I32").Value = 750
$ws.Range("J32").Value = 1761.25
$ws.Range("K32").Value = 750
$ws.Range("L32").Value = 1761.25
$ws.Range("M32").Value = -424
$ws.Range("N32").Value = -2413.25
$ws.Range("H33").Value = 368.37143
$ws.Range("I33").Value = 341.86667
$ws.Range("K33").Value = 341.86667
$ws.Range("M33").Value = -112.86667
$ws.Range("H132").Value = 7583255.5
$ws.Range("I132").Value = 14499724
$ws.Range("J132").Value = 8075.7617
$ws.Range("K132").Value = 43499172
$ws.Range("L132").Value = 24227.2851
$ws.Range("M132").Value = -43496642
$ws.Range("N132").Value = -29287.2851
$ws.Range("H137").Value = 1551.7391
$ws.Range("J137").Value = 1535.1428
$ws.Range("L137").Value = 4605.428400000001
$ws.Range("N137").Value = -9705.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1068.4
$ws.Range("I2").Value = 970.3333
$ws.Range("J2").Value = 1110.4286
$ws.Range("K2").Value = 970.3333
$ws.Range("L2").Value = 1110.4286
$ws.Range("M2").Value = -857.3333
$ws.Range("N2").Value = -1336.4286
$ws.Range("H45").Value = 1229.6666
$ws.Range("I45").Value = 1069.5
$ws.Range("K45").Value = 1069.5
$ws.Range("M45").Value = -692.5
$ws.Range("H63").Value = 2130.0833
$ws.Range("I63").Value = 2024.92
$ws.Range("J63").Value = 2369.0908
$ws.Range("K63").Value = 2024.92
$ws.Range("L63").Value = 2369.0908
$ws.Range("M63").Value = -1338.92
$ws.Range("N63").Value = -3741.0908
$ws.Range("H66").Value = 2130.0833
$ws.Range("I66").Value = 2024.92
$ws.Range("J66").Value = 2369.0908
$ws.Range("K66").Value = 10124.6
$ws.Range("L66").Value = 11845.454
$ws.Range("M66").Value = -6692.6
$ws.Range("N66").Value = -18709.454
$ws.Range("H74").Value = 1820.55
$ws.Range("I74").Value = 1523.9445
$ws.Range("K74").Value = 1523.9445
$ws.Range("M74").Value = -649.9445000000001
$ws.Range("H77").Value = 1820.55
$ws.Range("I77").Value = 1523.9445
$ws.Range("K77").Value = 7619.7225
$ws.Range("M77").Value = -3251.7225
$ws.Range("H110").Value = 1379.1
$ws.Range("I110").Value = 898.75
$ws.Range("J110").Value = 3300.5
$ws.Range("K110").Value = 898.75
$ws.Range("L110").Value = 3300.5
$ws.Range("M110").Value = 1146.25
$ws.Range("N110").Value = -7390.5
$ws.Range("H114").Value = 25832.334
$ws.Range("J114").Value = 25832.334
$ws.Range("L114").Value = 25832.334
$ws.Range("N114").Value = -34510.334
$ws.Range("H116").Value = 1068.4
$ws.Range("I116").Value = 970.3333
$ws.Range("J116").Value = 1110.4286
$ws.Range("K116").Value = 970.3333
$ws.Range("L116").Value = 1110.4286
$ws.Range("M116").Value = 1323.6667
$ws.Range("N116").Value = -5698.4286
$ws.Range("H122").Value = 1034.7
$ws.Range("I122").Value = 1034.7
$ws.Range("K122").Value = 3104.1
$ws.Range("M122").Value = -654.1000000000004

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1068.4
$ws.Range("I3").Value = 970.3333
$ws.Range("J3").Value = 1110.4286
$ws.Range("K3").Value = 970.3333
$ws.Range("L3").Value = 1110.4286
$ws.Range("M3").Value = -856.3333
$ws.Range("N3").Value = -1338.4286
$ws.Range("H107").Value = 1335.1177
$ws.Range("I107").Value = 1121
$ws.Range("K107").Value = 1121
$ws.Range("M107").Value = 799
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1199.7894
$ws.Range("J31").Value = 1487
$ws.Range("L31").Value = 1487
$ws.Range("N31").Value = -2077
$ws.Range("H34").Value = 1199.7894
$ws.Range("J34").Value = 1487
$ws.Range("L34").Value = 1487
$ws.Range("N34").Value = -1891
$ws.Range("H99").Value = 1631.8572
$ws.Range("J99").Value = 1590.75
$ws.Range("L99").Value = 1590.75
$ws.Range("N99").Value = -4586.75
$ws.Range("H126").Value = 1631.8572
$ws.Range("J126").Value = 1590.75
$ws.Range("L126").Value = 4772.25
$ws.Range("N126").Value = -9712.25
$ws.Range("H132").Value = 6571.5
$ws.Range("I132").Value = 12691.444
$ws.Range("J132").Value = 2334.6155
$ws.Range("K132").Value = 38074.33199999999
$ws.Range("L132").Value = 7003.8465
$ws.Range("M132").Value = -35544.33199999999
$ws.Range("N132").Value = -12063.8465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 644.4643
$ws.Range("I113").Value = 545
$ws.Range("J113").Value = 656.4
$ws.Range("K113").Value = 1635
$ws.Range("L113").Value = 1969.2
$ws.Range("M113").Value = 535
$ws.Range("N113").Value = -6309.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10488.8
$ws.Range("I5").Value = 626
$ws.Range("J5").Value = 14715.714
$ws.Range("K5").Value = 626
$ws.Range("L5").Value = 14715.714
$ws.Range("M5").Value = -514
$ws.Range("N5").Value = -14939.714
$ws.Range("H123").Value = 10125.2
$ws.Range("J123").Value = 10125.2
$ws.Range("L123").Value = 10125.2
$ws.Range("N123").Value = -15025.2
$ws.Range("H126").Value = 1959.8823
$ws.Range("I126").Value = 1821.2
$ws.Range("K126").Value = 5463.6
$ws.Range("M126").Value = -2993.6
$ws.Range("H132").Value = 5431.222
$ws.Range("I132").Value = 5554.7144
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 16664.1432
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -14134.1432
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 62535452
$ws.Range("I122").Value = 125050000
$ws.Range("K122").Value = 375150000
$ws.Range("M122").Value = -375147550
$ws.Range("H136").Value = 7892.533
$ws.Range("I136").Value = 9970.728
$ws.Range("K136").Value = 29912.184
$ws.Range("M136").Value = -27362.184

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2467.7778
$ws.Range("I132").Value = 2243.8572
$ws.Range("K132").Value = 6731.571599999999
$ws.Range("M132").Value = -4201.571599999999
$ws.Range("H136").Value = 802.9487
$ws.Range("I136").Value = 674.74194
$ws.Range("J136").Value = 1299.75
$ws.Range("K136").Value = 2024.22582
$ws.Range("L136").Value = 3899.25
$ws.Range("M136").Value = 525.7741799999999
$ws.Range("N136").Value = -8999.25
